$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks numeric while keeping it stored as text,
# exactly like the inline strings already present in this sheet (e.g. "59.990.51").
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $savedStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $savedStyle
}

$ws.Range("D2").Value = '60.003.23'
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("D3").Value = '2.655.68'
$ws.Range("E3").Value = '  +2.53%  '
$ws.Range("E4").Value = '  -0.02%  '
Set-TextValue "D5" '536.44'
$ws.Range("E5").Value = '  +1.31%  '
Set-TextValue "D6" '146.00'
$ws.Range("E6").Value = '  +4.49%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +1.54%  '
$ws.Range("D9").Value = '2.676.13'
$ws.Range("E9").Value = '  +2.80%  '
$ws.Range("E10").Value = '  +3.73%  '
$ws.Range("E11").Value = '  +2.46%  '
$ws.Range("E12").Value = '  +1.66%  '
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("D14").Value = '3.143.87'
$ws.Range("E14").Value = '  +2.91%  '
$ws.Range("D15").Value = '59.918.99'
$ws.Range("E15").Value = '  +1.50%  '
Set-TextValue "D16" '21.30'
$ws.Range("D17").Value = '2.643.81'
$ws.Range("E17").Value = '  +2.20%  '
$ws.Range("E18").Value = '  +1.59%  '
Set-TextValue "D19" '345.23'
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("E20").Value = '  +2.43%  '
$ws.Range("E21").Value = '  +2.08%  '
Set-TextValue "D22" '6.38'
$ws.Range("E22").Value = '  -0.50%  '
$ws.Range("E23").Value = '  +0.28%  '
Set-TextValue "D24" '67.51'
$ws.Range("E24").Value = '  +0.23%  '
Set-TextValue "D25" '0.416'
$ws.Range("E25").Value = '  +2.77%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("E27").Value = '  +0.01%  '
Set-TextValue "D28" '7.37'
$ws.Range("E28").Value = '  +3.02%  '
$ws.Range("E29").Value = '  +2.73%  '
Set-TextValue "D31" '1.66'
$ws.Range("E31").Value = '  +3.02%  '
Set-TextValue "D32" '5.92'
$ws.Range("E32").Value = '  +1.15%  '
Set-TextValue "D33" '19.14'
$ws.Range("E33").Value = '  +1.94%  '
Set-TextValue "D34" '150.40'
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("E35").Value = '  +2.05%  '
$ws.Range("E36").Value = '  +3.36%  '
$ws.Range("E37").Value = '  +0.57%  '
Set-TextValue "D38" '0.846'
$ws.Range("E38").Value = '  +2.12%  '
$ws.Range("E39").Value = '  +0.48%  '
Set-TextValue "D40" '291.61'
$ws.Range("E40").Value = '  +8.51%  '
Set-TextValue "D41" '3.61'
$ws.Range("E41").Value = '  +2.35%  '
Set-TextValue "D42" '0.998'
$ws.Range("E42").Value = '  +0.04%  '
Set-TextValue "D43" '0.606'
$ws.Range("E43").Value = '  +2.01%  '
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D44" '0.0544'
$ws.Range("E44").Value = '  +5.08%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue "D45" '10.74'
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D46" '0.0957'
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").Value = '1.978.95'
$ws.Range("E47").Value = '  +1.11%  '
$ws.Range("E48").Value = '  +2.71%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D49" '18.59'
$ws.Range("E49").Value = '  +2.32%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D50" '4.59'
$ws.Range("E50").Value = '  -0.49%  '
Set-TextValue "D51" '109.95'
$ws.Range("E51").Value = '  -1.62%  '
